$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 6 to make room for the new "depth" attribute,
# pushing temp/sal/biosat/O2_Ar_ratio/ncp/k down by one row.
$ws.Rows.Item(6).Insert()

# Row 2: utc_datetime -> datetime_utc (definition/class/format unchanged)
$ws.Range("A2").Value = "datetime_utc"

# Row 3: matlab_datetime -> datetime_utc_matlab
$ws.Range("A3").Value = "datetime_utc_matlab"

# Row 4: lat -> latitude
$ws.Range("A4").Value = "latitude"

# Row 5: lon -> longitude
$ws.Range("A5").Value = "longitude"

# Row 6 (new): depth attribute
$ws.Range("A6").Value = "depth"
$ws.Range("B6").Value = "Depth of sample below sea surface. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$ws.Range("C6").Value = "numeric"
$ws.Range("D6").Value = "meter"

# Row 10 (was O2_Ar_ratio at old row 9, now shifted to row 10):
# rename attribute and update its definition text
$ws.Range("A10").Value = "O2_Ar_ratio_corrected"
$ws.Range("B10").Value = "Oxygen-argon ratio of EIMS sample from underway corrected for air values"

# Update sheet dimension and selected cell to match the saved state
$ws.Range("B20").Select()
